$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.004.80"
$ws.Range("E2").Value = "  +2.08%  "

$ws.Range("D3").Value = "1.704.69"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D4").Value = "'1.001"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'316.00"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").Value = "'1.001"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "'0.3997"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  +2.05%  "

$ws.Range("D8").Value = "'0.4036"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  -0.69%  "

$ws.Range("D9").Value = "'1.472"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  -1.25%  "

$ws.Range("D10").Value = "'53.16"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  +0.59%  "

$ws.Range("D11").Value = "'1.001"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").Value = "'0.08811"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("D13").Value = "'26.06"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = "  -3.35%  "

$ws.Range("D14").Value = "'7.480"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001351"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.966"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").Value = "1.733.06"
$ws.Range("E17").Value = "  +2.59%  "

$ws.Range("D18").Value = "'96.02"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  -2.01%  "

$ws.Range("D19").Value = "'0.07196"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").Value = "'20.74"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").Value = "'7.316"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("D22").Value = "'1.000"
$ws.Cells.Item(22, 4).Style = "Normal"

$ws.Range("D23").Value = "'14.32"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").Value = "25.011.07"
$ws.Range("E24").Value = "  +2.09%  "

$ws.Range("D25").Value = "'2.409"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  +3.60%  "

$ws.Range("D26").Value = "'2.944"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  -2.78%  "

$ws.Range("D27").Value = "'23.57"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = "  +3.49%  "

$ws.Range("D28").Value = "'6.075"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  +12.71%  "

$ws.Range("D29").Value = "'162.86"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  -2.82%  "

$ws.Range("D30").Value = "'152.08"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  +5.40%  "

$ws.Range("D31").Value = "'8.402"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  -1.45%  "

$ws.Range("D32").Value = "'2.661"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  +20.05%  "

$ws.Range("D33").Value = "1.923.55"
$ws.Range("E33").Value = "  +2.54%  "

$ws.Range("D34").Value = "'0.08611"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  -1.90%  "

$ws.Range("D35").Value = "'0.03161"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  +3.16%  "

$ws.Range("D36").Value = "'1.051"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  +1.09%  "

$ws.Range("D37").Value = "'7.213"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("D38").Value = "'0.2908"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "  +3.74%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.09687"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "  +5.47%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'10.98"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").Value = "'0.8259"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  +2.91%  "

$ws.Range("D42").Value = "'14.02"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  -1.33%  "

$ws.Range("D43").Value = "'1.482"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("D44").Value = "'17.02"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  -2.60%  "

$ws.Range("D45").Value = "'2.687"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").Value = "'0.7387"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("D47").Value = "'0.09275"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "  +13.69%  "

$ws.Range("D48").Value = "'4.253"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("D49").Value = "'1.404"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "  -0.92%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").Value = "'139.84"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  -0.95%  "
